$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value looks numeric, so Excel
# preserves the exact string (leading/trailing zeros, multi-dot formatting)
# instead of silently coercing the text to a floating point number.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '69.217.45'
$ws.Range("E2").Value = '  +1.25%  '

$ws.Range("D3").Value = '3.887.77'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '529.58'
$ws.Range("E5").Value = '  +9.02%  '

$ws.Range("D6").Value = '144.04'
$ws.Range("E6").Value = '  -1.29%  '

$ws.Range("E7").Value = '  -1.91%  '

$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '0.718'
$ws.Range("E9").Value = '  -3.23%  '

$ws.Range("E10").Value = '  -5.52%  '

$ws.Range("D11").Value = '0.0000333'
$ws.Range("E11").Value = '  -5.89%  '

$ws.Range("D12").Value = '41.93'
$ws.Range("E12").Value = '  -2.60%  '

$ws.Range("D13").Value = '4.509.76'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("D14").Value = '10.22'
$ws.Range("E14").Value = '  -2.19%  '

$ws.Range("D15").Value = '3.880.71'
$ws.Range("E15").Value = '  -1.59%  '

$ws.Range("D16").Value = '14.01'
$ws.Range("E16").Value = '  -1.59%  '

$ws.Range("E17").Value = '  +6.77%  '

$ws.Range("E18").Value = '  -1.44%  '

$ws.Range("D19").Value = '20.34'
$ws.Range("E19").Value = '  +1.85%  '

$ws.Range("D20").Value = '69.180.24'
$ws.Range("E20").Value = '  +1.05%  '

$ws.Range("D21").Value = '423.64'
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = '3.40'
$ws.Range("E22").Value = '  -4.70%  '

$ws.Range("D23").Value = '14.18'
$ws.Range("E23").Value = '  -3.70%  '

$ws.Range("D24").Value = '87.53'
$ws.Range("E24").Value = '  -2.12%  '

$ws.Range("E25").Value = '  +7.97%  '

$ws.Range("D26").Value = '11.32'
$ws.Range("E26").Value = '  -8.02%  '

$ws.Range("D27").Value = '10.60'
$ws.Range("E27").Value = '  -3.34%  '

$ws.Range("D28").Value = '36.32'
$ws.Range("E28").Value = '  -2.85%  '

$ws.Range("D29").Value = '696.75'
$ws.Range("E29").Value = '  -3.09%  '

$ws.Range("E30").Value = '  -1.61%  '

$ws.Range("E31").Value = '  -3.34%  '

$ws.Range("D32").Value = '2.83'
$ws.Range("E32").Value = '  -3.09%  '

$ws.Range("D33").Value = '67.94'
$ws.Range("E33").Value = '  +10.09%  '

$ws.Range("D34").Value = '0.432'
$ws.Range("E34").Value = '  +6.98%  '

$ws.Range("D35").Value = '5.95'
$ws.Range("E35").Value = '  -2.01%  '

$ws.Range("D36").Value = '0.0₃0858'
$ws.Range("E36").Value = '  -4.48%  '

$ws.Range("D37").Value = '40.09'
$ws.Range("E37").Value = '  -1.50%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("E39").Value = '  -0.32%  '

$ws.Range("E40").Value = '  -0.11%  '

$ws.Range("D41").Value = '3.27'
$ws.Range("E41").Value = '  +5.27%  '

$ws.Range("E42").Value = '  +6.58%  '

$ws.Range("D43").Value = '0.0482'
$ws.Range("E43").Value = '  -3.07%  '

$ws.Range("D44").Value = '2.79'
$ws.Range("E44").Value = '  -6.72%  '

$ws.Range("E45").Value = '  +1.85%  '

$ws.Range("E46").Value = '  -1.38%  '

$ws.Range("D47").Value = '2.99'
$ws.Range("E47").Value = '  +6.42%  '

$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = '0.000272'
$ws.Range("E48").Value = '  +10.68%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.755.21'
$ws.Range("E49").Value = '  +14.86%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '144.73'
$ws.Range("E50").Value = '  +1.03%  '

$ws.Range("E51").Value = '  -2.81%  '
